$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 190.66667
$ws.Range("I2").Value = 204.8
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 204.8
$ws.Range("L2").Value = 120
$ws.Range("M2").Value = -91.80000000000001
$ws.Range("N2").Value = -346
$ws.Range("H28").Value = 3932.7273
$ws.Range("I28").Value = 160.16667
$ws.Range("K28").Value = 160.16667
$ws.Range("M28").Value = 324.83333
$ws.Range("H29").Value = 4457.143
$ws.Range("I29").Value = 375
$ws.Range("J29").Value = 9900
$ws.Range("K29").Value = 1125
$ws.Range("L29").Value = 29700
$ws.Range("M29").Value = -844
$ws.Range("N29").Value = -30262
$ws.Range("H38").Value = 695.6667
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("H58").Value = 1156.5385
$ws.Range("J58").Value = 2999.75
$ws.Range("L58").Value = 8999.25
$ws.Range("N58").Value = -9299.25
$ws.Range("H111").Value = 100
$ws.Range("J111").Value = 100
$ws.Range("L111").Value = 300
$ws.Range("N111").Value = -6434
$ws.Range("H112").Value = 4958.3335
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 5870
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 17610
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -19826
$ws.Range("H113").Value = 18800
$ws.Range("J113").Value = 3253
$ws.Range("L113").Value = 3253
$ws.Range("N113").Value = -9761
$ws.Range("H129").Value = 912.85
$ws.Range("J129").Value = 887.2593000000001
$ws.Range("L129").Value = 2661.7779
$ws.Range("N129").Value = -12661.7779
$ws.Range("H131").Value = 2319.5
$ws.Range("J131").Value = 4100.8
$ws.Range("L131").Value = 12302.4
$ws.Range("N131").Value = -22382.4
$ws.Range("H132").Value = 952.8919
$ws.Range("I132").Value = 948.7646999999999
$ws.Range("K132").Value = 2846.2941
$ws.Range("M132").Value = -316.2941000000001
$ws.Range("H137").Value = 1488.4546
$ws.Range("I137").Value = 1151.9412
$ws.Range("K137").Value = 3455.8236
$ws.Range("M137").Value = -905.8235999999997
$ws.Range("H138").Value = 1822.8438
$ws.Range("I138").Value = 1616.7
$ws.Range("J138").Value = 2166.4167
$ws.Range("K138").Value = 4850.1
$ws.Range("L138").Value = 6499.250100000001
$ws.Range("M138").Value = 289.8999999999996
$ws.Range("N138").Value = -16779.2501
$ws.Range("N38").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4103.451
$ws.Range("I32").Value = 2187.5789
$ws.Range("K32").Value = 2187.5789
$ws.Range("M32").Value = -1900.5789
$ws.Range("H45").Value = 1897.75
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1897.75
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1897.75
$ws.Range("N45").Value = -2651.75
$ws.Range("H61").Value = 3139.625
$ws.Range("I61").Value = 2286
$ws.Range("K61").Value = 2286
$ws.Range("M61").Value = -2074
$ws.Range("H63").Value = 7300
$ws.Range("I63").Value = 7300
$ws.Range("K63").Value = 7300
$ws.Range("M63").Value = -6614
$ws.Range("H66").Value = 7300
$ws.Range("I66").Value = 7300
$ws.Range("K66").Value = 36500
$ws.Range("M66").Value = -33068
$ws.Range("H74").Value = 701
$ws.Range("I74").Value = 719.7
$ws.Range("J74").Value = 514
$ws.Range("K74").Value = 719.7
$ws.Range("L74").Value = 514
$ws.Range("M74").Value = 154.3
$ws.Range("N74").Value = -2262
$ws.Range("H77").Value = 701
$ws.Range("I77").Value = 719.7
$ws.Range("J77").Value = 514
$ws.Range("K77").Value = 3598.5
$ws.Range("L77").Value = 2570
$ws.Range("M77").Value = 769.5
$ws.Range("N77").Value = -11306
$ws.Range("H92").Value = 49824.5
$ws.Range("J92").Value = 49824.5
$ws.Range("L92").Value = 49824.5
$ws.Range("N92").Value = -54816.5
$ws.Range("H97").Value = 1237.25
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 2000
$ws.Range("N97").Value = -2992
$ws.Range("H102").Value = 2145.75
$ws.Range("I102").Value = 1333.4
$ws.Range("J102").Value = 3499.6667
$ws.Range("K102").Value = 1333.4
$ws.Range("L102").Value = 3499.6667
$ws.Range("M102").Value = 288.5999999999999
$ws.Range("N102").Value = -6743.6667
$ws.Range("H132").Value = 1319.5483
$ws.Range("I132").Value = 1036
$ws.Range("K132").Value = 3108
$ws.Range("M132").Value = -578
$ws.Range("H136").Value = 3139.625
$ws.Range("I136").Value = 2286
$ws.Range("K136").Value = 6858
$ws.Range("M136").Value = -4308
$ws.Range("M45").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2124.7144
$ws.Range("I20").Value = 2094.6667
$ws.Range("J20").Value = 2305
$ws.Range("K20").Value = 2094.6667
$ws.Range("L20").Value = 2305
$ws.Range("M20").Value = -1847.6667
$ws.Range("N20").Value = -2799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 113.22222
$ws.Range("I7").Value = 139.83333
$ws.Range("K7").Value = 139.83333
$ws.Range("M7").Value = -26.83332999999999
$ws.Range("H31").Value = 2745.85
$ws.Range("I31").Value = 1556.4375
$ws.Range("K31").Value = 1556.4375
$ws.Range("M31").Value = -1261.4375
$ws.Range("H34").Value = 2745.85
$ws.Range("I34").Value = 1556.4375
$ws.Range("K34").Value = 1556.4375
$ws.Range("M34").Value = -1354.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("H75").Value = 1878.75
$ws.Range("I75").Value = 1500
$ws.Range("J75").Value = 2005
$ws.Range("K75").Value = 4500
$ws.Range("L75").Value = 6015
$ws.Range("M75").Value = -3502
$ws.Range("N75").Value = -8011
$ws.Range("H78").Value = 1878.75
$ws.Range("I78").Value = 1500
$ws.Range("J78").Value = 2005
$ws.Range("K78").Value = 13500
$ws.Range("L78").Value = 18045
$ws.Range("M78").Value = -8508
$ws.Range("N78").Value = -28029
$ws.Range("H116").Value = 166668670
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("H117").Value = 35714664
$ws.Range("I117").Value = 639
$ws.Range("K117").Value = 1917
$ws.Range("M117").Value = 1525
$ws.Range("H130").Value = 1702.5
$ws.Range("J130").Value = 1975
$ws.Range("L130").Value = 5925
$ws.Range("N130").Value = -15965
$ws.Range("H131").Value = 12873.396
$ws.Range("I131").Value = 473.8
$ws.Range("J131").Value = 14043.17
$ws.Range("K131").Value = 1421.4
$ws.Range("L131").Value = 42129.51
$ws.Range("M131").Value = 3618.6
$ws.Range("N131").Value = -52209.51
$ws.Range("M22").Value = ""
$ws.Range("M27").Value = ""
$ws.Range("M116").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 400
$ws.Range("I107").Value = 314.2857
$ws.Range("K107").Value = 314.2857
$ws.Range("M107").Value = 1605.7143
$ws.Range("H132").Value = 1329395.5
$ws.Range("I132").Value = 1924804.2
$ws.Range("J132").Value = 6265
$ws.Range("K132").Value = 5774412.6
$ws.Range("L132").Value = 18795
$ws.Range("M132").Value = -5771882.6
$ws.Range("N132").Value = -23855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5949.0625
$ws.Range("I7").Value = 2536.125
$ws.Range("J7").Value = 9362
$ws.Range("K7").Value = 2536.125
$ws.Range("L7").Value = 9362
$ws.Range("M7").Value = -2424.125
$ws.Range("N7").Value = -9586
$ws.Range("H61").Value = 2552.4443
$ws.Range("I61").Value = 2541.5
$ws.Range("J61").Value = 2574.3333
$ws.Range("K61").Value = 2541.5
$ws.Range("L61").Value = 2574.3333
$ws.Range("M61").Value = -2339.5
$ws.Range("N61").Value = -2978.3333
$ws.Range("H113").Value = 2552.4443
$ws.Range("I113").Value = 2541.5
$ws.Range("J113").Value = 2574.3333
$ws.Range("K113").Value = 2541.5
$ws.Range("L113").Value = 2574.3333
$ws.Range("M113").Value = -371.5
$ws.Range("N113").Value = -6914.3333
$ws.Range("H126").Value = 5949.0625
$ws.Range("I126").Value = 2536.125
$ws.Range("J126").Value = 9362
$ws.Range("K126").Value = 7608.375
$ws.Range("L126").Value = 28086
$ws.Range("M126").Value = -5138.375
$ws.Range("N126").Value = -33026
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 510.42105
$ws.Range("I107").Value = 351.66666
$ws.Range("J107").Value = 1105.75
$ws.Range("K107").Value = 1054.99998
$ws.Range("L107").Value = 3317.25
$ws.Range("M107").Value = 865.0000199999999
$ws.Range("N107").Value = -7157.25
$ws.Range("H136").Value = 24157430
$ws.Range("I136").Value = 32682500
$ws.Range("K136").Value = 98047500
$ws.Range("M136").Value = -98044950
